$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "4 Owned" -> "# Owned"
$ws.Range("C1").Value = "# Owned"

# Change values in C2:C5 from "Yes" (string) to 4 (number)
$ws.Range("C2:C5").Value = 4

# Move selection to I5 to match target
$ws.Range("I5").Select()
